$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 195
$ws.Range("A195").Value = "Government ERP UAT Lead / Manager (CGI Advantage 4.0)"
$ws.Range("B195").Value = "https://www.dice.com/job-detail/95bd6184-6a93-4f15-ad74-d3747bfd1ed9?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=5&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C195").Value = "New York"
$ws.Range("D195").Value = "Third Party, Contract"
$ws.Range("E195").Value = "USD 27.55 - 31.90 per hour"
$ws.Range("F195").Value = "PETADATA"

# New row 196
$ws.Range("A196").Value = "Google Cloud Platform Architect with AI"
$ws.Range("B196").Value = "https://www.dice.com/job-detail/938c6f1d-6b9f-473d-a3a6-11ae2db753ba?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=5&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C196").Value = "New York, New York"
$ws.Range("D196").Value = "Contract, Third Party"
$ws.Range("E196").Value = "Depends on Experience"
$ws.Range("F196").Value = "Tekfortune Inc."
